$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value = "v8_no_fl"
$ws.Range("B2").Value = "Versión 8 – sin fluorescencia"
$ws.Range("C2").Value = 0
$ws.Range("D2").Value = 3000
$ws.Range("E2").Value = 4000
$ws.Range("N2").Value = "0, 0.1, 1.0, 5.0, 10.0"
$ws.Range("R2").Value = "0, 0.004, 0.008, 0.012"
$ws.Range("W2").Value = "0.01, 0.018, 0.03"
